# Trade #57 closed at 2026-02-17 15:43:33 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.35   # Current Capital
$summary.Range("B4").Value = 0.35      # Total P&L $
$summary.Range("B5").Value = 0.12      # Total P&L %
$summary.Range("B6").Value = 57        # Total Trades
$summary.Range("B7").Value = 17        # Winning Trades
$summary.Range("B9").Value = 29.82     # Win Rate %

# ---- Strategy Status sheet ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.35     # Capital
$status.Range("D4").Value = 57         # Trades
$status.Range("E4").Value = 0.35       # P&L $
$status.Range("F4").Value = 0.35       # P&L %
$status.Range("G4").Value = 29.82      # Win Rate %

# ---- New trade row appended to "All Trades" and "MarketMaking" sheets ----
$newRow = 58

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item($newRow, 1).Value = 57
    $ws.Cells.Item($newRow, 2).Value = "2026-02-17"
    $ws.Cells.Item($newRow, 3).Value = "15:43:27"
    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"
    $ws.Cells.Item($newRow, 5).Value = "DOWN"
    $ws.Cells.Item($newRow, 6).Value = 0.8100000000000001
    $ws.Cells.Item($newRow, 7).Value = 0.86
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"
    $ws.Cells.Item($newRow, 9).Value = 6.1728
    $ws.Cells.Item($newRow, 10).Value = 0.05
    $ws.Cells.Item($newRow, 11).Value = 100.35
    $ws.Cells.Item($newRow, 12).Value = 0
    $ws.Cells.Item($newRow, 13).Value = 0
    $ws.Cells.Item($newRow, 14).Value = 0.6
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($newRow, 16).Value = "early_exit"
    $ws.Cells.Item($newRow, 17).Value = 0.14
}
